# Remove the "JUANCITO" worksheet, keeping only "PINI"
$excel.DisplayAlerts = $false
$wb = $excel.ActiveWorkbook

foreach ($sheet in @($wb.Worksheets)) {
    if ($sheet.Name -eq "JUANCITO") {
        $sheet.Delete()
    }
}

# Add the extra rows of data to the remaining "PINI" sheet
$ws = $wb.Worksheets.Item("PINI")

$ws.Range("A2").Value = "asdsa"
$ws.Range("B2").Value = "aadsa"
$ws.Range("A3").Value = "asdsa"
$ws.Range("B3").Value = "aadsa"
